$d = $word.ActiveDocument

# Merge "Versi" + "on" back into a single "Version" run (undo the earlier split).
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# Change " 2" back to " 1." (this merges with the following "." text in Word's model).
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2)

# Remove the now-superfluous trailing "." run that sits after the _GoBack bookmark.
$bm = $d.Bookmarks("_GoBack")
$tailStart = $bm.End
$p = $d.Paragraphs(1)
$tailEnd = $p.Range.End - 1
$tail = $d.Range($tailStart, $tailEnd)
$tail.Text = ""
